$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary totals: "VALOR MORA" total and "Cant. Periodos" count,
# since one of the two mora periods (2505) is being removed, leaving only 2506.
$ws.Range("E11").Value = 74314
$ws.Range("F13").Value = 1

# Remove the data row for period "2505" (row 17), keeping only the "2506" row.
# This shifts the trailing signature rows (22, 23) up to (21, 22).
$ws.Rows("17").Delete()
